$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell C2 (RunMode for TC001) from "Yes" to "No"
$ws.Range("C2").Value = "No"

# Row 1 - new headers (D1:J1)
$ws.Range("D1").Value = "Product Type"
$ws.Range("E1").Value = "Products to Add"
$ws.Range("F1").Value = "Sizes"
$ws.Range("G1").Value = "Colour"
$ws.Range("H1").Value = "QuantityToAdd"
$ws.Range("I1").Value = "Products to Remove"
$ws.Range("J1").Value = "QuantityToRemove"

# Row 3 - TC002 (A3:H3)
$ws.Range("A3").Value = "TC002"
$ws.Range("B3").Value = "Add to Cart"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "Women Bras&Tanks"
$ws.Range("E3").Value = "Celeste Sports Bra"
$ws.Range("F3").Value = "L"
$ws.Range("G3").Value = "Green"
$ws.Range("H3").Value = 2

# Row 4 - TC003 (A4:J4)
$ws.Range("A4").Value = "TC003"
$ws.Range("B4").Value = "Add and Remove from Cart"
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "Women Bras&Tanks"
$ws.Range("E4").Value = "Celeste Sports Bra"
$ws.Range("F4").Value = "L"
$ws.Range("G4").Value = "Green"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = "Celeste Sports Bra"
$ws.Range("J4").Value = 1

# Row 5 - TC004 (A5:H5)
$ws.Range("A5").Value = "TC004"
$ws.Range("B5").Value = "End to End Test 1"
$ws.Range("C5").Value = "No"
$ws.Range("D5").Value = "Men Jacket,Men Pants,Men Tanks,Women Tees"
$ws.Range("E5").Value = "Beaumont Summit Kit,Geo Insulated Jogging Pant,Rocco Gym Tank,Layla Tee"
$ws.Range("F5").Value = "L,34,M,S"
$ws.Range("G5").Value = "Red,Green,Blue,Red"
$ws.Range("H5").Value = "2,1,2,1"

# --- Formatting: reuse existing cell styles (avoid creating new style entries) ---
# Header row style (themed fill, centered) -> copy from A1 onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("D1:J1").PasteSpecial(-4122)

# Data row style (centered, no fill) -> copy from A2 onto the new populated data cells only
$ws.Range("A2").Copy()
$ws.Range("D3:H3").PasteSpecial(-4122)
$ws.Range("D4:J4").PasteSpecial(-4122)
$ws.Range("D5:H5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column widths (autofit results captured from the target file)
$ws.Columns.Item(2).ColumnWidth = 23.33203125
$ws.Columns.Item(4).ColumnWidth = 40.33203125
$ws.Columns.Item(5).ColumnWidth = 64.21875
$ws.Columns.Item(7).ColumnWidth = 17.21875
$ws.Columns.Item(8).ColumnWidth = 13.33203125
$ws.Columns.Item(9).ColumnWidth = 17.6640625
$ws.Columns.Item(10).ColumnWidth = 16.88671875

# View state: scroll so column C is the left-most visible column, select J4
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J4").Select()
